# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps on row 4 of the
# zh-cn and de-de language sheets to reflect the new handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-22 17:27:11"
$wsZhCn.Range("G4").Value = "2016-02-22 17:27:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-22 17:27:22"
$wsDeDe.Range("G4").Value = "2016-02-22 17:28:23"
